$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 228.18182
$ws.Range("I2").Value = 241
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 241
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = -128
$ws.Range("N2").Value = -326

$ws.Range("H106").Value = 62526.23
$ws.Range("I106").Value = 2733.0908
$ws.Range("J106").Value = 391388.5
$ws.Range("K106").Value = 2733.0908
$ws.Range("L106").Value = 391388.5
$ws.Range("M106").Value = -2102.0908
$ws.Range("N106").Value = -392650.5

$ws.Range("H135").Value = 20834058
$ws.Range("I135").Value = 579.35297
$ws.Range("J135").Value = 71429650
$ws.Range("K135").Value = 5214.17673
$ws.Range("L135").Value = 642866850
$ws.Range("M135").Value = -2679.17673
$ws.Range("N135").Value = -642871920

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 29999
$ws.Range("I109").Value = 29999
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 29999
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -28612
$ws.Range("N109").Value = ""

$ws.Range("H122").Value = 1737.2354
$ws.Range("I122").Value = 1663.4286
$ws.Range("J122").Value = 2081.6667
$ws.Range("K122").Value = 4990.2858
$ws.Range("L122").Value = 6245.000100000001
$ws.Range("M122").Value = -2540.2858
$ws.Range("N122").Value = -11145.0001

$ws.Range("H132").Value = 13516479
$ws.Range("I132").Value = 33335520
$ws.Range("J132").Value = 3496.7727
$ws.Range("K132").Value = 100006560
$ws.Range("L132").Value = 10490.3181
$ws.Range("M132").Value = -100004030
$ws.Range("N132").Value = -15550.3181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1133.3334
$ws.Range("I7").Value = 1133.3334
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1133.3334
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1020.3334
$ws.Range("N7").Value = ""

$ws.Range("H30").Value = 10011
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 10011
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 10011
$ws.Range("M30").Value = ""
$ws.Range("N30").Value = -10261

$ws.Range("H86").Value = 3588.7646
$ws.Range("I86").Value = 3187.3333
$ws.Range("J86").Value = 6599.5
$ws.Range("K86").Value = 3187.3333
$ws.Range("L86").Value = 6599.5
$ws.Range("M86").Value = -2064.3333
$ws.Range("N86").Value = -8845.5

$ws.Range("H89").Value = 3588.7646
$ws.Range("I89").Value = 3187.3333
$ws.Range("J89").Value = 6599.5
$ws.Range("K89").Value = 15936.6665
$ws.Range("L89").Value = 32997.5
$ws.Range("M89").Value = -10320.6665
$ws.Range("N89").Value = -44229.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 105818.55
$ws.Range("I6").Value = 505000
$ws.Range("J6").Value = 17111.555
$ws.Range("K6").Value = 505000
$ws.Range("L6").Value = 17111.555
$ws.Range("M6").Value = -504887
$ws.Range("N6").Value = -17337.555

$ws.Range("H22").Value = 1513.625
$ws.Range("I22").Value = 202.54546
$ws.Range("J22").Value = 4398
$ws.Range("K22").Value = 202.54546
$ws.Range("L22").Value = 4398
$ws.Range("M22").Value = 147.45454
$ws.Range("N22").Value = -5098

$ws.Range("H28").Value = 40000
$ws.Range("J28").Value = 40000
$ws.Range("L28").Value = 40000
$ws.Range("N28").Value = -40490

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 5074.7827
$ws.Range("I113").Value = 6775.125
$ws.Range("J113").Value = 1188.2858
$ws.Range("K113").Value = 20325.375
$ws.Range("L113").Value = 3564.8574
$ws.Range("M113").Value = -18155.375
$ws.Range("N113").Value = -7904.857400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 2008106.9
$ws.Range("I2").Value = 2008106.9
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2008106.9
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -2007993.9
$ws.Range("N2").Value = ""

$ws.Range("H31").Value = 10332.75
$ws.Range("I31").Value = 1110.3334
$ws.Range("J31").Value = 38000
$ws.Range("K31").Value = 1110.3334
$ws.Range("L31").Value = 38000
$ws.Range("M31").Value = -818.3334
$ws.Range("N31").Value = -38584

$ws.Range("H37").Value = 10332.75
$ws.Range("I37").Value = 1110.3334
$ws.Range("J37").Value = 38000
$ws.Range("K37").Value = 1110.3334
$ws.Range("L37").Value = 38000
$ws.Range("M37").Value = -833.3334
$ws.Range("N37").Value = -38554

$ws.Range("H122").Value = 2895.7778
$ws.Range("I122").Value = 3007.75
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 9023.25
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -6573.25
$ws.Range("N122").Value = -10900

$ws.Range("H132").Value = 4078.4
$ws.Range("I132").Value = 2799.5
$ws.Range("J132").Value = 5996.75
$ws.Range("K132").Value = 8398.5
$ws.Range("L132").Value = 17990.25
$ws.Range("M132").Value = -5868.5
$ws.Range("N132").Value = -23050.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1153.6923
$ws.Range("I22").Value = 1066.6666
$ws.Range("J22").Value = 1179.8
$ws.Range("K22").Value = 1066.6666
$ws.Range("L22").Value = 1179.8
$ws.Range("M22").Value = -771.6666
$ws.Range("N22").Value = -1769.8

$ws.Range("H27").Value = 1153.6923
$ws.Range("I27").Value = 1066.6666
$ws.Range("J27").Value = 1179.8
$ws.Range("K27").Value = 1066.6666
$ws.Range("L27").Value = 1179.8
$ws.Range("M27").Value = -959.6666
$ws.Range("N27").Value = -1393.8

$ws.Range("H108").Value = 41996.5
$ws.Range("J108").Value = 41996.5
$ws.Range("L108").Value = 41996.5
$ws.Range("N108").Value = -49676.5

$ws.Range("H123").Value = 28208.572
$ws.Range("J123").Value = 28208.572
$ws.Range("L123").Value = 28208.572
$ws.Range("N123").Value = -38008.572

$ws.Range("H133").Value = 36000
$ws.Range("J133").Value = 36000
$ws.Range("L133").Value = 36000
$ws.Range("N133").Value = -41060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = ""

$ws.Range("H109").Value = 44377
$ws.Range("J109").Value = 44377
$ws.Range("L109").Value = 44377
$ws.Range("N109").Value = -47151

$ws.Range("H123").Value = 34871.668
$ws.Range("J123").Value = 34871.668
$ws.Range("L123").Value = 34871.668
$ws.Range("N123").Value = -44671.668

$ws.Range("H133").Value = 27500
$ws.Range("J133").Value = 27500
$ws.Range("L133").Value = 27500
$ws.Range("N133").Value = -37620
